$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 591.6667
$ws.Range("J19").Value = 591.6667
$ws.Range("L19").Value = 591.6667
$ws.Range("N19").Value = -941.6667

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5250
$ws.Range("I43").Value = 5666
$ws.Range("J43").Value = 4002
$ws.Range("K43").Value = 5666
$ws.Range("L43").Value = 4002
$ws.Range("M43").Value = -5597
$ws.Range("N43").Value = -4140

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 12000.968
$ws.Range("I74").Value = 12780.909
$ws.Range("J74").Value = 10094.444
$ws.Range("K74").Value = 12780.909
$ws.Range("L74").Value = 10094.444
$ws.Range("M74").Value = -11844.909
$ws.Range("N74").Value = -11966.444

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 12000.968
$ws.Range("I77").Value = 12780.909
$ws.Range("J77").Value = 10094.444
$ws.Range("K77").Value = 63904.545
$ws.Range("L77").Value = 50472.22
$ws.Range("M77").Value = -59224.545
$ws.Range("N77").Value = -59832.22

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 200002290
$ws.Range("I86").Value = 250002160
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 250002160
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -250001037
$ws.Range("N86").Value = -5046

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 200002290
$ws.Range("I89").Value = 250002160
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 1250010800
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -1250005184
$ws.Range("N89").Value = -25232

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 4079.4
$ws.Range("I103").Value = 5669.3335
$ws.Range("K103").Value = 17008.0005
$ws.Range("M103").Value = -16422.0005

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3295.3333
$ws.Range("I113").Value = 3376.8
$ws.Range("K113").Value = 3376.8
$ws.Range("M113").Value = -122.8000000000002

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 98389.5
$ws.Range("J123").Value = 98389.5
$ws.Range("L123").Value = 98389.5
$ws.Range("N123").Value = -108189.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9079.666999999999
$ws.Range("I132").Value = 10605.375
$ws.Range("K132").Value = 31816.125
$ws.Range("M132").Value = -29286.125

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1356901.1
$ws.Range("I137").Value = 4547188
$ws.Range("J137").Value = 7164.423
$ws.Range("K137").Value = 13641564
$ws.Range("L137").Value = 21493.269
$ws.Range("M137").Value = -13639014
$ws.Range("N137").Value = -26593.269

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3430.323
$ws.Range("I138").Value = 2743.25
$ws.Range("J138").Value = 3735.689
$ws.Range("K138").Value = 8229.75
$ws.Range("L138").Value = 11207.067
$ws.Range("M138").Value = -3089.75
$ws.Range("N138").Value = -21487.067

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2671.1345
$ws.Range("I32").Value = 2671.1345
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2671.1345
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2384.1345
$ws.Range("N32").ClearContents()

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 30396.4
$ws.Range("I45").Value = 39448.363
$ws.Range("J45").Value = 5503.5
$ws.Range("K45").Value = 39448.363
$ws.Range("L45").Value = 5503.5
$ws.Range("M45").Value = -39071.363
$ws.Range("N45").Value = -6257.5

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2568.5217
$ws.Range("I61").Value = 2230.0454
$ws.Range("K61").Value = 2230.0454
$ws.Range("M61").Value = -2018.0454

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3189.394
$ws.Range("I132").Value = 2728.4783
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 8185.4349
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -5655.4349
$ws.Range("N132").Value = -17808.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2568.5217
$ws.Range("I136").Value = 2230.0454
$ws.Range("K136").Value = 6690.1362
$ws.Range("M136").Value = -4140.1362

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1749.0952
$ws.Range("J86").Value = 2020.1
$ws.Range("L86").Value = 2020.1
$ws.Range("N86").Value = -4266.1

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1749.0952
$ws.Range("J89").Value = 2020.1
$ws.Range("L89").Value = 10100.5
$ws.Range("N89").Value = -21332.5

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13686541
$ws.Range("J105").Value = 41670116
$ws.Range("L105").Value = 41670116
$ws.Range("N105").Value = -41673610

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3475
$ws.Range("I134").Value = 3229.8276
$ws.Range("J134").Value = 4363.75
$ws.Range("K134").Value = 9689.4828
$ws.Range("L134").Value = 13091.25
$ws.Range("M134").Value = -7154.4828
$ws.Range("N134").Value = -18161.25

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3253.4
$ws.Range("J31").Value = 6005.2666
$ws.Range("L31").Value = 6005.2666
$ws.Range("N31").Value = -6595.2666

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3253.4
$ws.Range("J34").Value = 6005.2666
$ws.Range("L34").Value = 6005.2666
$ws.Range("N34").Value = -6409.2666

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 71432584
$ws.Range("I132").Value = 125001780
$ws.Range("K132").Value = 375005340
$ws.Range("M132").Value = -375002810

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4063
$ws.Range("I134").Value = 4170.4287
$ws.Range("K134").Value = 12511.2861
$ws.Range("M134").Value = -9976.286100000001

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 67078576
$ws.Range("I4").Value = 78343096
$ws.Range("K4").Value = 235029288
$ws.Range("M4").Value = -235029176

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4141.0713
$ws.Range("J81").Value = 4679.5454
$ws.Range("L81").Value = 14038.6362
$ws.Range("N81").Value = -16284.6362

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 4141.0713
$ws.Range("J84").Value = 4679.5454
$ws.Range("L84").Value = 42115.9086
$ws.Range("N84").Value = -53347.9086

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 240045.8
$ws.Range("I121").Value = 33410
$ws.Range("J121").Value = 549999.5
$ws.Range("K121").Value = 100230
$ws.Range("L121").Value = 1649998.5
$ws.Range("M121").Value = -98920
$ws.Range("N121").Value = -1652618.5

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2999
$ws.Range("I133").Value = 2999
$ws.Range("K133").Value = 8997
$ws.Range("M133").Value = -3937

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58826264
$ws.Range("I80").Value = 100002216
$ws.Range("K80").Value = 100002216
$ws.Range("M80").Value = -100001218

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 58826264
$ws.Range("I83").Value = 100002216
$ws.Range("K83").Value = 500011080
$ws.Range("M83").Value = -500006088

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8633.357
$ws.Range("I102").Value = 1499.5714
$ws.Range("K102").Value = 1499.5714
$ws.Range("M102").Value = 122.4286

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3587
$ws.Range("I122").Value = 2403.2632
$ws.Range("J122").Value = 5836.1
$ws.Range("K122").Value = 7209.7896
$ws.Range("L122").Value = 17508.3
$ws.Range("M122").Value = -4759.7896
$ws.Range("N122").Value = -22408.3

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2940.95
$ws.Range("J132").Value = 4102.4
$ws.Range("L132").Value = 12307.2
$ws.Range("N132").Value = -17367.2

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4622.8057
$ws.Range("I40").Value = 4948.84
$ws.Range("J40").Value = 3881.818
$ws.Range("K40").Value = 4948.84
$ws.Range("L40").Value = 3881.818
$ws.Range("M40").Value = -4812.84
$ws.Range("N40").Value = -4153.818

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 14356.056
$ws.Range("I61").Value = 2529
$ws.Range("J61").Value = 55750.75
$ws.Range("K61").Value = 2529
$ws.Range("L61").Value = 55750.75
$ws.Range("M61").Value = -2327
$ws.Range("N61").Value = -56154.75

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1299.8
$ws.Range("I93").Value = 1166.6666
$ws.Range("J93").Value = 1499.5
$ws.Range("K93").Value = 1166.6666
$ws.Range("L93").Value = 1499.5
$ws.Range("M93").Value = 81.33339999999998
$ws.Range("N93").Value = -3995.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 14356.056
$ws.Range("I113").Value = 2529
$ws.Range("J113").Value = 55750.75
$ws.Range("K113").Value = 2529
$ws.Range("L113").Value = 55750.75
$ws.Range("M113").Value = -359
$ws.Range("N113").Value = -60090.75

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5888.3
$ws.Range("I132").Value = 3419.5557
$ws.Range("J132").Value = 7908.1816
$ws.Range("K132").Value = 10258.6671
$ws.Range("L132").Value = 23724.5448
$ws.Range("M132").Value = -7728.667099999999
$ws.Range("N132").Value = -28784.5448

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 20835690
$ws.Range("I122").Value = 2559.8
$ws.Range("J122").Value = 35716496
$ws.Range("K122").Value = 7679.400000000001
$ws.Range("L122").Value = 107149488
$ws.Range("M122").Value = -5229.400000000001
$ws.Range("N122").Value = -107154388

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3205.6667
$ws.Range("I126").Value = 3309.5
$ws.Range("J126").Value = 2998
$ws.Range("K126").Value = 9928.5
$ws.Range("L126").Value = 8994
$ws.Range("M126").Value = -7458.5
$ws.Range("N126").Value = -13934

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8335733.5
$ws.Range("I132").Value = 10102926
$ws.Range("K132").Value = 30308778
$ws.Range("M132").Value = -30306248

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10063.538
$ws.Range("I136").Value = 9961.333000000001
$ws.Range("J136").Value = 10293.5
$ws.Range("K136").Value = 29883.999
$ws.Range("L136").Value = 30880.5
$ws.Range("M136").Value = -27333.999
$ws.Range("N136").Value = -35980.5
